# Add two new review rows (rows 6 and 7) to Sheet1, mirroring the
# formatting of the existing rows and wiring up mailto: hyperlinks for
# the email columns (C and D), exactly like the existing rows 2-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 (com.singleton.strechy / taxi game review) -------------------
# Copy formatting from row 2, which already uses the same per-column
# style pattern (A=1, B=0, C=2, D=2, E=0, F=1) that row 6 needs.
$ws.Range("A2:F2").Copy()
$ws.Range("A6:F6").PasteSpecial(-4122)

$ws.Range("A6").Value = "com.singleton.strechy"
$ws.Range("B6").Value = "taxi game"
$ws.Range("C6").Value = "redvelvetmichael@gmail.com"
$ws.Range("D6").Value = "veredsnir12@gmail.com"
$ws.Range("E6").Value = "27/5/2019 15:59"
$ws.Range("F6").Value = "Crazy and hard levels but I like it. I can play it all day long. Sweet taxi alabama!!"

$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:redvelvetmichael@gmail.com", "", "", "redvelvetmichael@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:veredsnir12@gmail.com", "", "", "veredsnir12@gmail.com")

# Hyperlinks.Add() stamps its own built-in "Hyperlink" look onto the
# cell, clobbering the style we just copied in. Re-apply the plain
# email-column formatting (from row 2) on top so C6/D6 stay visually
# consistent with the rest of the table, exactly like C5/D5 already do.
$ws.Range("C2:D2").Copy()
$ws.Range("C6:D6").PasteSpecial(-4122)

# --- Row 7 (com.hamxa.shaynachim / bitcoin guide review) -----------------
# Copy formatting from row 5, which already uses the same per-column
# style pattern (A=1, B=0, C=2, D=2, E=0, F=1) that row 7 needs.
$ws.Range("A5:F5").Copy()
$ws.Range("A7:F7").PasteSpecial(-4122)

$ws.Range("A7").Value = "com.hamxa.shaynachim"
$ws.Range("B7").Value = "bitcoin guide"
$ws.Range("C7").Value = "vikicrestina@gmail.com"
$ws.Range("D7").Value = "cristianjohn1222@gmail.com"
$ws.Range("E7").Value = "27/5/2019 15:59"
$ws.Range("F7").Value = "bitcoin guide – great app. Following KISS guidelines – Keep it simple st…"

$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:vikicrestina@gmail.com", "", "", "vikicrestina@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:cristianjohn1222@gmail.com", "", "", "cristianjohn1222@gmail.com")

# Same style fix-up as row 6 above.
$ws.Range("C5:D5").Copy()
$ws.Range("C7:D7").PasteSpecial(-4122)

# Keep the active selection on F7, matching the source workbook's cursor
# position after the new rows were appended.
$ws.Range("F7").Select()
